$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("G10").Value = 3
$ws.Range("I10").Value = 2.27
$ws.Range("J10").Value = 3.45
$ws.Range("K10").Value = 2.07
$ws.Range("L10").Value = 2.85
$ws.Range("M10").Value = 1.33
$ws.Range("N10").Value = 2.8
$ws.Range("O10").Value = 1.98
$ws.Range("P10").Value = 1.65
$ws.Range("Q10").Value = 3.2
$ws.Range("R10").Value = 1.26
$ws.Range("U10").Value = 1.75
$ws.Range("V10").Value = 1.87
$ws.Range("X10").Value = 15.5
$ws.Range("Y10").Value = 10.75
$ws.Range("Z10").Value = 37
$ws.Range("AA10").Value = 27
$ws.Range("AB10").Value = 35
$ws.Range("AC10").Value = 8.75
$ws.Range("AD10").Value = 6.1
$ws.Range("AF10").Value = 70
$ws.Range("AG10").Value = 600
$ws.Range("AH10").Value = 7.3
$ws.Range("AJ10").Value = 9
$ws.Range("AK10").Value = 23
$ws.Range("AL10").Value = 19.5
$ws.Range("AM10").Value = 30

# Row 11
$ws.Range("O11").Value = 1.67
$ws.Range("P11").Value = 2.15
$ws.Range("AN11").Value = 1.03
$ws.Range("AO11").Value = 10

# Row 12
$ws.Range("G12").Value = 2.2
$ws.Range("I12").Value = 2.9
$ws.Range("J12").Value = 2.75
$ws.Range("L12").Value = 3.25
$ws.Range("Y12").Value = 9.5
$ws.Range("Z12").Value = 21
$ws.Range("AB12").Value = 21
$ws.Range("AJ12").Value = 12

# Row 13
$ws.Range("G13").Value = 2.25
$ws.Range("H13").Value = 3.1
$ws.Range("J13").Value = 3
$ws.Range("K13").Value = 2.05
$ws.Range("M13").Value = 1.33
$ws.Range("N13").Value = 3.25
$ws.Range("O13").Value = 2.05
$ws.Range("P13").Value = 1.75
$ws.Range("Q13").Value = 3.5
$ws.Range("R13").Value = 1.29
$ws.Range("S13").Value = 1.44
$ws.Range("T13").Value = 2.63
$ws.Range("U13").Value = 1.83
$ws.Range("V13").Value = 1.83
$ws.Range("W13").Value = 8
$ws.Range("Y13").Value = 10
$ws.Range("AB13").Value = 29
$ws.Range("AC13").Value = 9
$ws.Range("AD13").Value = 6
$ws.Range("AE13").Value = 15
$ws.Range("AF13").Value = 51
$ws.Range("AG13").Value = 700
$ws.Range("AH13").Value = 9.5
$ws.Range("AL13").Value = 26
$ws.Range("AN13").Value = 1.04
$ws.Range("AO13").Value = 9

# Row 17
$ws.Range("G17").Value = 4.2
$ws.Range("I17").Value = 1.85
$ws.Range("J17").Value = 4.75
$ws.Range("L17").Value = 2.5
$ws.Range("M17").Value = 1.3
$ws.Range("N17").Value = 3.4
$ws.Range("O17").Value = 2
$ws.Range("P17").Value = 1.8
$ws.Range("Q17").Value = 3.5
$ws.Range("R17").Value = 1.29
$ws.Range("Y17").Value = 15
$ws.Range("AI17").Value = 8.5
$ws.Range("AJ17").Value = 8.5
$ws.Range("AK17").Value = 15
$ws.Range("AL17").Value = 15

